$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 3).Value = 53250
$ws.Cells.Item(9, 3).Value = 49350
$ws.Cells.Item(12, 3).Value = 45870
$ws.Cells.Item(13, 3).Value = 1200
$ws.Cells.Item(19, 3).Value = 5700
$ws.Cells.Item(21, 3).Value = 9080
$ws.Cells.Item(25, 3).Value = 44950
$ws.Cells.Item(26, 3).Value = 84450
$ws.Cells.Item(37, 3).Value = 20230
$ws.Cells.Item(38, 3).Value = 47110
$ws.Cells.Item(40, 3).Value = 60170
$ws.Cells.Item(44, 3).Value = 29300
$ws.Cells.Item(45, 3).Value = 21130
$ws.Cells.Item(47, 3).Value = 48300
$ws.Cells.Item(48, 3).Value = 41090
$ws.Cells.Item(50, 3).Value = 41790
$ws.Cells.Item(51, 3).Value = 4450
$ws.Cells.Item(56, 3).Value = 6650
$ws.Cells.Item(57, 3).Value = 19760
$ws.Cells.Item(60, 3).Value = 2320
$ws.Cells.Item(63, 3).Value = 14760
$ws.Cells.Item(66, 3).Value = 61210
$ws.Cells.Item(68, 3).Value = 68120
$ws.Cells.Item(69, 3).Value = 40930
$ws.Cells.Item(70, 3).Value = 33740
$ws.Cells.Item(74, 3).Value = 8080
$ws.Cells.Item(78, 3).Value = 30620
$ws.Cells.Item(86, 3).Value = 17360
$ws.Cells.Item(87, 3).Value = 70840
$ws.Cells.Item(88, 3).Value = 16500
$ws.Cells.Item(90, 3).Value = 3900
$ws.Cells.Item(105, 3).Value = 2020
$ws.Cells.Item(107, 3).Value = 80640
$ws.Cells.Item(109, 3).Value = 41020
$ws.Cells.Item(117, 3).Value = 21980
$ws.Cells.Item(130, 3).Value = 24500
$ws.Cells.Item(131, 3).Value = 55540
$ws.Cells.Item(144, 3).Value = 63170
